# Auto-generated Excel COM-interop script
# Applies numeric value updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 875.06665
$ws.Range("I28").Value = 823.9167
$ws.Range("K28").Value = 823.9167
$ws.Range("M28").Value = -338.9167
$ws.Range("H137").Value = 65710.5
$ws.Range("I137").Value = 127424.25
$ws.Range("J137").Value = 3996.75
$ws.Range("K137").Value = 382272.75
$ws.Range("L137").Value = 11990.25
$ws.Range("M137").Value = -379722.75
$ws.Range("N137").Value = -17090.25
$ws.Range("H138").Value = 2828.1516
$ws.Range("I138").Value = 2344.3333
$ws.Range("J138").Value = 4118.3335
$ws.Range("K138").Value = 7032.999899999999
$ws.Range("L138").Value = 12355.0005
$ws.Range("M138").Value = -1892.999899999999
$ws.Range("N138").Value = -22635.0005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17302.912
$ws.Range("I32").Value = 16932.277
$ws.Range("J32").Value = 25333.334
$ws.Range("K32").Value = 16932.277
$ws.Range("L32").Value = 25333.334
$ws.Range("M32").Value = -16645.277
$ws.Range("N32").Value = -25907.334
$ws.Range("H61").Value = 12407.595
$ws.Range("I61").Value = 13539.866
$ws.Range("J61").Value = 7555
$ws.Range("K61").Value = 13539.866
$ws.Range("L61").Value = 7555
$ws.Range("M61").Value = -13327.866
$ws.Range("N61").Value = -7979
$ws.Range("H74").Value = 25611.373
$ws.Range("I74").Value = 27284.525
$ws.Range("K74").Value = 27284.525
$ws.Range("M74").Value = -26410.525
$ws.Range("H77").Value = 25611.373
$ws.Range("I77").Value = 27284.525
$ws.Range("K77").Value = 136422.625
$ws.Range("M77").Value = -132054.625
$ws.Range("H122").Value = 3928.3333
$ws.Range("I122").Value = 3916
$ws.Range("K122").Value = 11748
$ws.Range("M122").Value = -9298
$ws.Range("H132").Value = 25954.357
$ws.Range("I132").Value = 29692.527
$ws.Range("K132").Value = 89077.58099999999
$ws.Range("M132").Value = -86547.58099999999
$ws.Range("H136").Value = 12407.595
$ws.Range("I136").Value = 13539.866
$ws.Range("J136").Value = 7555
$ws.Range("K136").Value = 40619.598
$ws.Range("L136").Value = 22665
$ws.Range("M136").Value = -38069.598
$ws.Range("N136").Value = -27765

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3079.3125
$ws.Range("I86").Value = 2618
$ws.Range("K86").Value = 2618
$ws.Range("M86").Value = -1495
$ws.Range("H89").Value = 3079.3125
$ws.Range("I89").Value = 2618
$ws.Range("K89").Value = 13090
$ws.Range("M89").Value = -7474
$ws.Range("H134").Value = 1954.5
$ws.Range("I134").Value = 1773.375
$ws.Range("J134").Value = 5577
$ws.Range("K134").Value = 5320.125
$ws.Range("L134").Value = 16731
$ws.Range("M134").Value = -2785.125
$ws.Range("N134").Value = -21801

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2559.4285
$ws.Range("I31").Value = 2277.32
$ws.Range("J31").Value = 4910.3335
$ws.Range("K31").Value = 2277.32
$ws.Range("L31").Value = 4910.3335
$ws.Range("M31").Value = -1982.32
$ws.Range("N31").Value = -5500.3335
$ws.Range("H34").Value = 2559.4285
$ws.Range("I34").Value = 2277.32
$ws.Range("J34").Value = 4910.3335
$ws.Range("K34").Value = 2277.32
$ws.Range("L34").Value = 4910.3335
$ws.Range("M34").Value = -2075.32
$ws.Range("N34").Value = -5314.3335
$ws.Range("H58").Value = 52171.65
$ws.Range("I58").Value = 57699.332
$ws.Range("K58").Value = 57699.332
$ws.Range("M58").Value = -57496.332
$ws.Range("H62").Value = 68399.664
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 68399.664
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H99").Value = 13162.556
$ws.Range("I99").Value = 34566
$ws.Range("K99").Value = 34566
$ws.Range("M99").Value = -33068
$ws.Range("H122").Value = 2672.1538
$ws.Range("I122").Value = 2031.125
$ws.Range("J122").Value = 3697.8
$ws.Range("K122").Value = 6093.375
$ws.Range("L122").Value = 11093.4
$ws.Range("M122").Value = -3643.375
$ws.Range("N122").Value = -15993.4
$ws.Range("H126").Value = 13162.556
$ws.Range("I126").Value = 34566
$ws.Range("K126").Value = 103698
$ws.Range("M126").Value = -101228
$ws.Range("H136").Value = 52171.65
$ws.Range("I136").Value = 57699.332
$ws.Range("K136").Value = 173097.996
$ws.Range("M136").Value = -170547.996

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1657689.8
$ws.Range("I4").Value = 1598002.4
$ws.Range("K4").Value = 4794007.199999999
$ws.Range("M4").Value = -4793895.199999999
$ws.Range("H13").Value = 74.40000000000001
$ws.Range("I13").Value = 74.40000000000001
$ws.Range("K13").Value = 223.2
$ws.Range("M13").Value = -55.20000000000002
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H68").Value = 11739.866
$ws.Range("H71").Value = 11739.866
$ws.Range("H103").Value = 1554.5294
$ws.Range("I103").Value = 634.7778
$ws.Range("J103").Value = 2589.25
$ws.Range("K103").Value = 1904.3334
$ws.Range("L103").Value = 7767.75
$ws.Range("M103").Value = -1025.3334
$ws.Range("N103").Value = -9525.75
$ws.Range("H120").Value = 18381.2
$ws.Range("I120").Value = 15476.75
$ws.Range("J120").Value = 29999
$ws.Range("K120").Value = 46430.25
$ws.Range("L120").Value = 89997
$ws.Range("M120").Value = -41592.25
$ws.Range("N120").Value = -99673
$ws.Range("H131").Value = 12508304
$ws.Range("J131").Value = 25015000
$ws.Range("L131").Value = 75045000
$ws.Range("N131").Value = -75055080

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 55.083332
$ws.Range("I2").Value = 56.25
$ws.Range("J2").Value = 52.75
$ws.Range("K2").Value = 56.25
$ws.Range("L2").Value = 52.75
$ws.Range("M2").Value = 56.75
$ws.Range("N2").Value = -278.75
$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42216
$ws.Range("H117").Value = 1500
$ws.Range("J117").Value = 1500
$ws.Range("L117").Value = 1500
$ws.Range("N117").Value = -8384
$ws.Range("H122").Value = 3249.5454
$ws.Range("I122").Value = 2677.8572
$ws.Range("K122").Value = 8033.571599999999
$ws.Range("M122").Value = -5583.571599999999
$ws.Range("H126").Value = 6109.6665
$ws.Range("I126").Value = 5495
$ws.Range("K126").Value = 16485
$ws.Range("M126").Value = -14015

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 139468.62
$ws.Range("I22").Value = 139468.62
$ws.Range("K22").Value = 139468.62
$ws.Range("M22").Value = -139173.62
$ws.Range("H27").Value = 139468.62
$ws.Range("I27").Value = 139468.62
$ws.Range("K27").Value = 139468.62
$ws.Range("M27").Value = -139361.62
$ws.Range("H61").Value = 3759.8125
$ws.Range("I61").Value = 2665.923
$ws.Range("K61").Value = 2665.923
$ws.Range("M61").Value = -2463.923
$ws.Range("H113").Value = 3759.8125
$ws.Range("I113").Value = 2665.923
$ws.Range("K113").Value = 2665.923
$ws.Range("M113").Value = -495.9229999999998
$ws.Range("H132").Value = 301874.75
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 99999
$ws.Range("J133").Value = 99999
$ws.Range("L133").Value = 99999
$ws.Range("N133").Value = -105059
$ws.Range("H136").Value = 3204.6155
$ws.Range("I136").Value = 2896.5757
$ws.Range("K136").Value = 8689.7271
$ws.Range("M136").Value = -6139.7271

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 18999.5
$ws.Range("J37").Value = 18999.5
$ws.Range("L37").Value = 18999.5
$ws.Range("N37").Value = -19405.5
$ws.Range("H113").Value = 911.4314000000001
$ws.Range("I113").Value = 708.91113
$ws.Range("K113").Value = 2126.73339
$ws.Range("M113").Value = 43.26661000000013
$ws.Range("H126").Value = 129385.94
$ws.Range("I126").Value = 162116.17
$ws.Range("K126").Value = 486348.51
$ws.Range("M126").Value = -483878.51
$ws.Range("H132").Value = 64675.766
$ws.Range("I132").Value = 64675.766
$ws.Range("K132").Value = 194027.298
$ws.Range("M132").Value = -191497.298
$ws.Range("H136").Value = 2763.804
$ws.Range("I136").Value = 2249.0417
$ws.Range("J136").Value = 11000
$ws.Range("K136").Value = 6747.125100000001
$ws.Range("L136").Value = 33000
$ws.Range("M136").Value = -4197.125100000001
$ws.Range("N136").Value = -38100

